$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wsAbout = $wb.Worksheets.Item("About")

# A2: standalone version string
$wsAbout.Range("A2").Value = "Version: $newVersion"

# A6: citation text containing the version string embedded
$a6 = [string]$wsAbout.Range("A6").Value2
$a6New = $a6 -replace [regex]::Escape($oldVersion), $newVersion
$wsAbout.Range("A6").Value = $a6New

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# S2:S9 contain the standalone version string (build_version column)
for ($r = 2; $r -le 9; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    if ([string]$cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
